# Auto-generated script to apply 2023-11-30 daily crime data update
# Updates column J (year 2023 totals) across Citywide Totals, By Neighborhood,
# and individual neighborhood sheets to reflect newly added data for 2023-11-30.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 10).Value = 7031
$ws.Cells.Item(3, 10).Value = 7420
$ws.Cells.Item(4, 10).Value = 1620
$ws.Cells.Item(5, 10).Value = 582
$ws.Cells.Item(6, 10).Value = 10029
$ws.Cells.Item(7, 10).Value = 26682

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(2, 10).Value = 70
$ws.Cells.Item(7, 10).Value = 393

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 10).Value = 445
$ws.Cells.Item(3, 10).Value = 496
$ws.Cells.Item(6, 10).Value = 602
$ws.Cells.Item(7, 10).Value = 1674

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(2, 10).Value = 155
$ws.Cells.Item(3, 10).Value = 200
$ws.Cells.Item(7, 10).Value = 535

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(4, 10).Value = 56
$ws.Cells.Item(5, 10).Value = 51
$ws.Cells.Item(6, 10).Value = 426
$ws.Cells.Item(7, 10).Value = 1206

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 10).Value = 245
$ws.Cells.Item(7, 10).Value = 824

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(2, 10).Value = 112
$ws.Cells.Item(7, 10).Value = 408

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 10).Value = 211
$ws.Cells.Item(6, 10).Value = 203
$ws.Cells.Item(7, 10).Value = 763
$ws.Cells.Item(8, 10).Value = 1674
$ws.Cells.Item(18, 10).Value = 220
$ws.Cells.Item(19, 10).Value = 772
$ws.Cells.Item(24, 10).Value = 86
$ws.Cells.Item(27, 10).Value = 161
$ws.Cells.Item(29, 10).Value = 1430
$ws.Cells.Item(33, 10).Value = 1206
$ws.Cells.Item(34, 10).Value = 122
$ws.Cells.Item(37, 10).Value = 824
$ws.Cells.Item(42, 10).Value = 1151
$ws.Cells.Item(44, 10).Value = 204
$ws.Cells.Item(48, 10).Value = 302
$ws.Cells.Item(50, 10).Value = 159
$ws.Cells.Item(53, 10).Value = 393
$ws.Cells.Item(54, 10).Value = 522
$ws.Cells.Item(55, 10).Value = 420
$ws.Cells.Item(57, 10).Value = 125
$ws.Cells.Item(63, 10).Value = 88
$ws.Cells.Item(67, 10).Value = 996
$ws.Cells.Item(77, 10).Value = 184
$ws.Cells.Item(78, 10).Value = 310
$ws.Cells.Item(80, 10).Value = 46
$ws.Cells.Item(83, 10).Value = 535
$ws.Cells.Item(85, 10).Value = 1104
$ws.Cells.Item(88, 10).Value = 286
$ws.Cells.Item(89, 10).Value = 334
$ws.Cells.Item(91, 10).Value = 308
$ws.Cells.Item(93, 10).Value = 112
$ws.Cells.Item(97, 10).Value = 243
$ws.Cells.Item(99, 10).Value = 408
$ws.Cells.Item(101, 10).Value = 26682

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(3, 10).Value = 372
$ws.Cells.Item(6, 10).Value = 276
$ws.Cells.Item(7, 10).Value = 996

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(2, 10).Value = 129
$ws.Cells.Item(4, 10).Value = 40
$ws.Cells.Item(6, 10).Value = 244
$ws.Cells.Item(7, 10).Value = 522

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 10).Value = 432
$ws.Cells.Item(7, 10).Value = 1430

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(3, 10).Value = 56
$ws.Cells.Item(7, 10).Value = 302

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(5, 10).Value = 30
$ws.Cells.Item(6, 10).Value = 298
$ws.Cells.Item(7, 10).Value = 772

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(6, 10).Value = 81
$ws.Cells.Item(7, 10).Value = 204

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Cells.Item(3, 10).Value = 48
$ws.Cells.Item(7, 10).Value = 203

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(3, 10).Value = 230
$ws.Cells.Item(6, 10).Value = 613
$ws.Cells.Item(7, 10).Value = 1151

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(2, 10).Value = 82
$ws.Cells.Item(7, 10).Value = 310

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(6, 10).Value = 237
$ws.Cells.Item(7, 10).Value = 420

$ws = $wb.Worksheets.Item('Dunning')
$ws.Cells.Item(3, 10).Value = 22
$ws.Cells.Item(7, 10).Value = 86

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(6, 10).Value = 79
$ws.Cells.Item(7, 10).Value = 308

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Cells.Item(3, 10).Value = 45
$ws.Cells.Item(7, 10).Value = 220

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Cells.Item(2, 10).Value = 32
$ws.Cells.Item(7, 10).Value = 112

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(3, 10).Value = 231
$ws.Cells.Item(5, 10).Value = 20
$ws.Cells.Item(7, 10).Value = 763

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Cells.Item(3, 10).Value = 32
$ws.Cells.Item(7, 10).Value = 122

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Cells.Item(4, 10).Value = 24
$ws.Cells.Item(7, 10).Value = 159

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(2, 10).Value = 63
$ws.Cells.Item(7, 10).Value = 211

$ws = $wb.Worksheets.Item('West Town')
$ws.Cells.Item(6, 10).Value = 168
$ws.Cells.Item(7, 10).Value = 243

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(2, 10).Value = 58
$ws.Cells.Item(7, 10).Value = 286

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(3, 10).Value = 97
$ws.Cells.Item(7, 10).Value = 334

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(2, 10).Value = 42
$ws.Cells.Item(7, 10).Value = 161

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Cells.Item(6, 10).Value = 58
$ws.Cells.Item(7, 10).Value = 125

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 10).Value = 296
$ws.Cells.Item(3, 10).Value = 397
$ws.Cells.Item(7, 10).Value = 1104

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Cells.Item(3, 10).Value = 60
$ws.Cells.Item(4, 10).Value = 17
$ws.Cells.Item(7, 10).Value = 184

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Cells.Item(6, 10).Value = 24
$ws.Cells.Item(7, 10).Value = 46
